$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, $rowA, $rowB) {
    $rangeA = $sheet.Range("B$rowA`:AD$rowA")
    $rangeB = $sheet.Range("B$rowB`:AD$rowB")
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

# Simple pairwise swaps of B:AD content (column A -- the running index -- stays put)
Swap-Rows $ws 10 11
Swap-Rows $ws 38 39
Swap-Rows $ws 47 48
Swap-Rows $ws 63 64
Swap-Rows $ws 68 69
Swap-Rows $ws 81 82
Swap-Rows $ws 86 87
Swap-Rows $ws 123 124
Swap-Rows $ws 147 148

# 3-way rotation among rows 35, 36, 37: new_35 = old_36, new_36 = old_37, new_37 = old_35
$range35 = $ws.Range("B35:AD35")
$range36 = $ws.Range("B36:AD36")
$range37 = $ws.Range("B37:AD37")

$val35 = $range35.Value2
$val36 = $range36.Value2
$val37 = $range37.Value2

$range35.Value2 = $val36
$range36.Value2 = $val37
$range37.Value2 = $val35

Write-Output "done"
